$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear L41 (was "WIP") -> now blank, task becomes DONE
$ws.Range("L41").Value = $null
$ws.Range("L41").NumberFormat = "General"

# Fill in the new task row 47: "LightSensor Beacon"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "LightSensor Beacon"

$ws.Range("H47").NumberFormat = "@"
$ws.Range("H47").Value = "LightSensor code"

$ws.Range("K47").NumberFormat = "@"
$ws.Range("K47").Value = "LB"
